# Applies scheduled-runner profit/price updates to the Diabolos_Profits sheets.
# Values correspond to refreshed Universalis market data pulled for each
# crafting-class leve (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38

$ws.Range("H39").Value = 35714664
$ws.Range("I39").Value = 304.6
$ws.Range("J39").Value = 125000560
$ws.Range("K39").Value = 913.8000000000001
$ws.Range("L39").Value = 375001680
$ws.Range("M39").Value = -617.8000000000001
$ws.Range("N39").Value = -375002272

$ws.Range("H132").Value = 3732.5625
$ws.Range("I132").Value = 3378.2856
$ws.Range("J132").Value = 6212.5
$ws.Range("K132").Value = 10134.8568
$ws.Range("L132").Value = 18637.5
$ws.Range("M132").Value = -7604.856800000001
$ws.Range("N132").Value = -23697.5

$ws.Range("H138").Value = 2561.353
$ws.Range("I138").Value = 1285.2778
$ws.Range("J138").Value = 3996.9375
$ws.Range("K138").Value = 3855.8334
$ws.Range("L138").Value = 11990.8125
$ws.Range("M138").Value = 1284.1666
$ws.Range("N138").Value = -22270.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2462.386
$ws.Range("I32").Value = 2421.9058
$ws.Range("K32").Value = 2421.9058
$ws.Range("M32").Value = -2134.9058

$ws.Range("H61").Value = 55558510
$ws.Range("I61").Value = 66669572
$ws.Range("K61").Value = 66669572
$ws.Range("M61").Value = -66669360

$ws.Range("H63").Value = 80007730
$ws.Range("I63").Value = 125007240
$ws.Range("J63").Value = 28579714
$ws.Range("K63").Value = 125007240
$ws.Range("L63").Value = 28579714
$ws.Range("M63").Value = -125006554
$ws.Range("N63").Value = -28581086

$ws.Range("H66").Value = 80007730
$ws.Range("I66").Value = 125007240
$ws.Range("J66").Value = 28579714
$ws.Range("K66").Value = 625036200
$ws.Range("L66").Value = 142898570
$ws.Range("M66").Value = -625032768
$ws.Range("N66").Value = -142905434

$ws.Range("H132").Value = 83335940
$ws.Range("I132").Value = 90911570
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 272734710
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -272732180
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 55558510
$ws.Range("I136").Value = 66669572
$ws.Range("K136").Value = 200008716
$ws.Range("M136").Value = -200006166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1667.1915
$ws.Range("I86").Value = 1538.1034
$ws.Range("J86").Value = 1875.1666
$ws.Range("K86").Value = 1538.1034
$ws.Range("L86").Value = 1875.1666
$ws.Range("M86").Value = -415.1034
$ws.Range("N86").Value = -4121.1666

$ws.Range("H89").Value = 1667.1915
$ws.Range("I89").Value = 1538.1034
$ws.Range("J89").Value = 1875.1666
$ws.Range("K89").Value = 7690.517
$ws.Range("L89").Value = 9375.833000000001
$ws.Range("M89").Value = -2074.517
$ws.Range("N89").Value = -20607.833

$ws.Range("H94").Value = 908.5172
$ws.Range("I94").Value = 719.5454999999999
$ws.Range("J94").Value = 1502.4286
$ws.Range("K94").Value = 719.5454999999999
$ws.Range("L94").Value = 1502.4286
$ws.Range("M94").Value = -268.5454999999999
$ws.Range("N94").Value = -2404.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1635.4
$ws.Range("I22").Value = 799.5714
$ws.Range("J22").Value = 2366.75
$ws.Range("K22").Value = 799.5714
$ws.Range("L22").Value = 2366.75
$ws.Range("M22").Value = -449.5714
$ws.Range("N22").Value = -3066.75

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 34933.223
$ws.Range("J74").Value = 34933.223
$ws.Range("L74").Value = 34933.223
$ws.Range("N74").Value = -36681.223

$ws.Range("H77").Value = 34933.223
$ws.Range("J77").Value = 34933.223
$ws.Range("L77").Value = 104799.669
$ws.Range("N77").Value = -113535.669

$ws.Range("H107").Value = 1718.9524
$ws.Range("I107").Value = 1849.8948
$ws.Range("K107").Value = 1849.8948
$ws.Range("M107").Value = 70.10519999999997

$ws.Range("H122").Value = 1642.8096
$ws.Range("I122").Value = 1699
$ws.Range("J122").Value = 519
$ws.Range("K122").Value = 5097
$ws.Range("L122").Value = 1557
$ws.Range("M122").Value = -2647
$ws.Range("N122").Value = -6457

$ws.Range("H132").Value = 1728.1
$ws.Range("I132").Value = 1450.6842
$ws.Range("K132").Value = 4352.0526
$ws.Range("M132").Value = -1822.0526

$ws.Range("H134").Value = 3244.1667
$ws.Range("I134").Value = 2693.2
$ws.Range("K134").Value = 8079.599999999999
$ws.Range("M134").Value = -5544.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1275.75
$ws.Range("I5").Value = 804
$ws.Range("J5").Value = 1747.5
$ws.Range("K5").Value = 2412
$ws.Range("L5").Value = 5242.5
$ws.Range("M5").Value = -2300
$ws.Range("N5").Value = -5466.5

$ws.Range("H7").Value = 33.333332
$ws.Range("I7").Value = 37
$ws.Range("K7").Value = 111
$ws.Range("M7").Value = 1

$ws.Range("H26").Value = 187
$ws.Range("I26").Value = 80.5
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 241.5
$ws.Range("L26").Value = 1200
$ws.Range("M26").Value = 46.5
$ws.Range("N26").Value = -1776

$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -3996

$ws.Range("H135").Value = 1275.75
$ws.Range("I135").Value = 804
$ws.Range("J135").Value = 1747.5
$ws.Range("K135").Value = 7236
$ws.Range("L135").Value = 15727.5
$ws.Range("M135").Value = -4701
$ws.Range("N135").Value = -20797.5

$ws.Range("H140").Value = 1801.9
$ws.Range("I140").Value = 1186.5
$ws.Range("J140").Value = 2725
$ws.Range("K140").Value = 3559.5
$ws.Range("L140").Value = 8175
$ws.Range("M140").Value = 1620.5
$ws.Range("N140").Value = -18535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 128.1875
$ws.Range("I2").Value = 85.888885
$ws.Range("K2").Value = 85.888885
$ws.Range("M2").Value = 27.111115

$ws.Range("H15").Value = 57999
$ws.Range("J15").Value = 57999
$ws.Range("L15").Value = 57999
$ws.Range("N15").Value = -58575

$ws.Range("H81").Value = 57999
$ws.Range("J81").Value = 57999
$ws.Range("L81").Value = 57999
$ws.Range("N81").Value = -59995

$ws.Range("H84").Value = 57999
$ws.Range("J84").Value = 57999
$ws.Range("L84").Value = 173997
$ws.Range("N84").Value = -183981

$ws.Range("H113").Value = 3323.818
$ws.Range("J113").Value = 4045.375
$ws.Range("L113").Value = 4045.375
$ws.Range("N113").Value = -8385.375

$ws.Range("H126").Value = 8689.632
$ws.Range("I126").Value = 12561.3
$ws.Range("J126").Value = 4387.778
$ws.Range("K126").Value = 37683.89999999999
$ws.Range("L126").Value = 13163.334
$ws.Range("M126").Value = -35213.89999999999
$ws.Range("N126").Value = -18103.334

$ws.Range("H132").Value = 3886.3215
$ws.Range("I132").Value = 3548.8
$ws.Range("K132").Value = 10646.4
$ws.Range("M132").Value = -8116.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23813024
$ws.Range("J7").Value = 4250.4
$ws.Range("L7").Value = 4250.4
$ws.Range("N7").Value = -4474.4

$ws.Range("H55").Value = 237.8
$ws.Range("I55").Value = 237.8
$ws.Range("K55").Value = 237.8
$ws.Range("M55").Value = -64.80000000000001

$ws.Range("H61").Value = 5859.44
$ws.Range("I61").Value = 3468.4211
$ws.Range("J61").Value = 13431
$ws.Range("K61").Value = 3468.4211
$ws.Range("L61").Value = 13431
$ws.Range("M61").Value = -3266.4211
$ws.Range("N61").Value = -13835

$ws.Range("H113").Value = 5859.44
$ws.Range("I113").Value = 3468.4211
$ws.Range("J113").Value = 13431
$ws.Range("K113").Value = 3468.4211
$ws.Range("L113").Value = 13431
$ws.Range("M113").Value = -1298.4211
$ws.Range("N113").Value = -17771

$ws.Range("H122").Value = 3386.64
$ws.Range("I122").Value = 2666.6316
$ws.Range("K122").Value = 7999.8948
$ws.Range("M122").Value = -5549.8948

$ws.Range("H126").Value = 23813024
$ws.Range("J126").Value = 4250.4
$ws.Range("L126").Value = 12751.2
$ws.Range("N126").Value = -17691.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9099061
$ws.Range("J81").Value = 13342229
$ws.Range("L81").Value = 26684458
$ws.Range("N81").Value = -26686580

$ws.Range("H84").Value = 9099061
$ws.Range("J84").Value = 13342229
$ws.Range("L84").Value = 133422290
$ws.Range("N84").Value = -133432898

$ws.Range("H132").Value = 4769.76
$ws.Range("I132").Value = 4560.263
$ws.Range("J132").Value = 5433.1665
$ws.Range("K132").Value = 13680.789
$ws.Range("L132").Value = 16299.4995
$ws.Range("M132").Value = -11150.789
$ws.Range("N132").Value = -21359.4995

$ws.Range("H136").Value = 3637.5789
$ws.Range("I136").Value = 1646.8182
$ws.Range("J136").Value = 6374.875
$ws.Range("K136").Value = 4940.4546
$ws.Range("L136").Value = 19124.625
$ws.Range("M136").Value = -2390.4546
$ws.Range("N136").Value = -24224.625
